$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 330, pushing the existing rows 330-379 down to 331-380.
$ws.Rows("330").Insert()

# Populate the newly inserted row 330 with the new record.
$ws.Cells.Item(330, 1).Value = 4
$ws.Cells.Item(330, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(330, 3).Value = "Los Lagos"
$ws.Cells.Item(330, 4).Value = 44984
$ws.Cells.Item(330, 5).Value = 10
$ws.Cells.Item(330, 6).Value = "Fruta"
$ws.Cells.Item(330, 7).Value = 100108
$ws.Cells.Item(330, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(330, 9).Value = 100108005
$ws.Cells.Item(330, 10).Value = "Piña"
$ws.Cells.Item(330, 11).Value = "Caramelo"
$ws.Cells.Item(330, 12).Value = "Primera"
$ws.Cells.Item(330, 13).Value = 120
$ws.Cells.Item(330, 14).Value = 26000
$ws.Cells.Item(330, 15).Value = 27000
$ws.Cells.Item(330, 16).Value = 26500
$ws.Cells.Item(330, 17).Value = "$/caja 12 unidades"
$ws.Cells.Item(330, 18).Value = "Ecuador"
$ws.Cells.Item(330, 19).Value = 2208
$ws.Cells.Item(330, 20).Value = 12
